# Auto-update price data: insert a new "today" row at the top (row 2),
# pushing all existing data rows down by one, mirroring the commit
# "自动更新价格数据 2026-01-07 02:51:52".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand new row right below the header row.
$ws.Rows.Item(2).Insert()
# Excel's Insert() copies formatting from the row above (the bold header);
# strip that back out so the new row matches the plain data-row style.
$ws.Rows.Item(2).ClearFormats()

# Write the new date as literal text (not an Excel date serial number).
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value() = "2026-01-07"
# Remove the now-unneeded "@" text format so the cell keeps the default style.
$ws.Range("A2").ClearFormats()

$ws.Range("B2").Value() = 783.5
$ws.Range("C2").Value() = 1112
$ws.Range("D2").Value() = 3610
